# Add new data row (row 20) for "BIOTA SYNTHESIS" / Sao Paulo, Brazil.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting the row (rather than just writing into an empty row) makes Excel
# carry the formatting down from row 19 into row 20 for us, matching the
# style pattern already used throughout the table.
$ws.Rows("20").Insert()

$ws.Range("A20").Value = "BIOTA SYNTHESIS"
$ws.Range("B20").Value = "BIOTA SYNTHESIS"
# Column B on this row is unlike the rows above it (which keep the
# inherited style) - it matches column A's unstyled look instead.
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "São Paulo"
$ws.Range("D20").Value = "São Paulo"
$ws.Range("E20").Value = "Brazil"
$ws.Range("F20").Value = "BR"
$ws.Range("G20").Value = -23.561140000000002
$ws.Range("H20").Value = -46.722847000000002
$ws.Range("I20").Value = "Yes"

# Move / record the active selection like the saved workbook shows.
$ws.Range("G21").Select()
